# commiting hide row and dates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

# Add the new "originalIndex" column header in V1, matching the style of
# the other header cells (e.g. A1) by cloning formats from A1.
$ws.Range("A1").Copy()
$ws.Range("V1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("V1").Value = "originalIndex"

# Fill V2:V51 with the original (0-based) row index before any reordering.
for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 22).Value = $row - 2
}

# Fill in the previously-blank Billing Start Date / Remarks for the last row.
# Force the NumberFormat to Text first so Excel doesn't auto-convert the
# "12/25/2023" literal into a date serial, then clear the format so the
# cell is left on the default (unstyled) format, matching the other blank
# text cells on this row.
$ws.Range("R51").NumberFormat = "@"
$ws.Range("R51").Value = "12/25/2023"
$ws.Range("R51").ClearFormats()

$ws.Range("T51").Value = "-"
